# Weekly update: a new daily price record is inserted for Femacal de La
# Calera - Papa right after row 305 (new row 306), pushing every
# subsequent record down by one row (old row 306 -> new row 307, ...,
# old row 354 -> new row 355).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 306; this shifts rows 306..354 down to 307..355
# and carries the column-D date style (s="2") onto the new row, same as Excel's
# native "Insert" behaviour.
$ws.Rows.Item(306).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(306, 1).Value  = 3
$ws.Cells.Item(306, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(306, 3).Value  = "Coquimbo"
$ws.Cells.Item(306, 4).Value  = 44505
$ws.Cells.Item(306, 5).Value  = 5
$ws.Cells.Item(306, 6).Value  = 100114001
$ws.Cells.Item(306, 7).Value  = "Papa"
$ws.Cells.Item(306, 8).Value  = "Asterix"
$ws.Cells.Item(306, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(306, 10).Value = 530
$ws.Cells.Item(306, 11).Value = 10000
$ws.Cells.Item(306, 12).Value = 11000
$ws.Cells.Item(306, 13).Value = 10472
$ws.Cells.Item(306, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(306, 15).Value = "Provincia de Talca"
$ws.Cells.Item(306, 16).Value = 419
$ws.Cells.Item(306, 17).Value = 25
$ws.Cells.Item(306, 18).Value = "Hortaliza"
